$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price/Volume columns so values like "307.50" or "113.00"
# are preserved exactly as text rather than being auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.196.95"
$ws.Range("E2").Value = "  -2.50%  "
$ws.Range("D3").Value = "1.873.89"
$ws.Range("E3").Value = "  -1.89%  "
$ws.Range("D4").Value = "0.9990"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "307.50"
$ws.Range("E5").Value = "  -1.64%  "
$ws.Range("D6").Value = "0.9997"
$ws.Range("D7").Value = "0.5125"
$ws.Range("E7").Value = "  +2.59%  "
$ws.Range("D8").Value = "0.3749"
$ws.Range("E8").Value = "  -1.36%  "
$ws.Range("D9").Value = "0.07143"
$ws.Range("E9").Value = "  -1.90%  "
$ws.Range("D10").Value = "0.8882"
$ws.Range("E10").Value = "  -2.41%  "
$ws.Range("D11").Value = "20.67"
$ws.Range("E11").Value = "  -2.86%  "
$ws.Range("D12").Value = "1.883.44"
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").Value = "0.07547"
$ws.Range("E13").Value = "  -1.20%  "
$ws.Range("D14").Value = "5.327"
$ws.Range("E14").Value = "  -2.63%  "
$ws.Range("D15").Value = "89.28"
$ws.Range("E15").Value = "  -3.47%  "
$ws.Range("D16").Value = "0.9993"
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").Value = "0.000008482"
$ws.Range("E17").Value = "  -2.77%  "
$ws.Range("D18").Value = "14.09"
$ws.Range("E18").Value = "  -3.84%  "
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").Value = "27.223.56"
$ws.Range("E20").Value = "  -2.56%  "
$ws.Range("D21").Value = "5.057"
$ws.Range("E21").Value = "  -2.17%  "
$ws.Range("D22").Value = "2.115.67"
$ws.Range("E22").Value = "  -2.44%  "
$ws.Range("D23").Value = "10.57"
$ws.Range("E23").Value = "  -2.81%  "
$ws.Range("D24").Value = "6.479"
$ws.Range("E24").Value = "  -1.94%  "
$ws.Range("D25").Value = "150.07"
$ws.Range("E25").Value = "  -1.62%  "
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("D27").Value = "17.98"
$ws.Range("E27").Value = "  -2.22%  "
$ws.Range("D28").Value = "2.101"
$ws.Range("E28").Value = "  -5.31%  "
$ws.Range("D29").Value = "113.00"
$ws.Range("E29").Value = "  -1.85%  "
$ws.Range("D30").Value = "4.725"
$ws.Range("E30").Value = "  -3.41%  "
$ws.Range("D31").Value = "4.675"
$ws.Range("E31").Value = "  -3.15%  "
$ws.Range("D32").Value = "0.09034"
$ws.Range("E32").Value = "  +0.61%  "
$ws.Range("D33").Value = "0.05146"
$ws.Range("E33").Value = "  -2.65%  "
$ws.Range("E34").Value = "  -2.39%  "
$ws.Range("D35").Value = "1.159"
$ws.Range("E35").Value = "  -6.46%  "
$ws.Range("D36").Value = "0.7360"
$ws.Range("E36").Value = "  -6.44%  "
$ws.Range("D37").Value = "0.02055"
$ws.Range("E37").Value = "  -1.29%  "
$ws.Range("D38").Value = "2.514"
$ws.Range("E38").Value = "  -5.62%  "
$ws.Range("D39").Value = "3.061"
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("E40").Value = "  -1.45%  "
$ws.Range("D41").Value = "0.5375"
$ws.Range("E41").Value = "  -2.74%  "
$ws.Range("D42").Value = "6.590"
$ws.Range("E42").Value = "  -3.16%  "
$ws.Range("D43").Value = "117.27"
$ws.Range("E43").Value = "  +2.98%  "
$ws.Range("D44").Value = "8.351"
$ws.Range("E44").Value = "  -1.94%  "
$ws.Range("D45").Value = "0.1474"
$ws.Range("E45").Value = "  -2.41%  "
$ws.Range("D46").Value = "0.4642"
$ws.Range("E46").Value = "  -3.67%  "
$ws.Range("D47").Value = "0.9998"
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("D48").Value = "10.07"
$ws.Range("E48").Value = "  -4.99%  "
$ws.Range("D49").Value = "1.573"
$ws.Range("E49").Value = "  -3.91%  "
$ws.Range("D50").Value = "64.46"
$ws.Range("E50").Value = "  -4.45%  "
$ws.Range("D51").Value = "36.54"
$ws.Range("E51").Value = "  -1.12%  "
